$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The author deleted the row for "The Space Paradox" / "November 28"
# (row 73), which shifts every following row up by one (135 rows -> 134).
$ws.Rows.Item(73).Delete()

# Re-apply the sheet's sort (on column B, the data is/was already sorted
# by date-within-group) over the now-smaller range so the persisted
# sortState/sortCondition refs shrink from A1:B135/B1:B135 to
# A1:B134/B1:B134 to match the new data extent.
$sortRange = $ws.Range("A1:B134")
$keyRange = $ws.Range("B1:B134")
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($keyRange)
$ws.Sort.SetRange($sortRange)
$ws.Sort.Header = -4135
$ws.Sort.Apply()

# Restore the author's final viewport/selection position after the edit.
$null = $ws.Range("B127").Select()
$excel.ActiveWindow.ScrollRow = 113
$excel.ActiveWindow.ScrollColumn = 1
